$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 10295.35
$ws.Range("B9").Value = 10379.42
$ws.Range("C9").Value = 107.89
$ws.Range("D9").Value = 107.02
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -0.81
$ws.Range("G9").Value = 42613.765509259261
$ws.Range("H9").Value = $false

# Row 10
$ws.Range("A10").Value = 10254.17
$ws.Range("B10").Value = 10295.35
$ws.Range("C10").Value = 107.17
$ws.Range("D10").Value = 106.74
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = -0.4
$ws.Range("G10").Value = 42614.672789351855
$ws.Range("H10").Value = $false

# Row 11
$ws.Range("A11").Value = 10189.57
$ws.Range("B11").Value = 10254.17
$ws.Range("C11").Value = 107.04
$ws.Range("D11").Value = 106.37
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = -0.63
$ws.Range("G11").Value = 42615.750115740739
$ws.Range("H11").Value = $false
